$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 3700
$ws.Range("I2").Value = 9613
$ws.Range("J2").Value = 40001
$ws.Range("K2").Value = 172
$ws.Range("L2").Value = 10861
$ws.Range("M2").Value = 667
$ws.Range("N2").Value = 7109
$ws.Range("O2").Value = 28
$ws.Range("P2").Value = 159
$ws.Range("Q2").Value = 60
$ws.Range("R2").Value = 532
$ws.Range("S2").Value = 4252
$ws.Range("T2").Value = 7028
$ws.Range("U2").Value = 525
$ws.Range("V2").Value = 62089
$ws.Range("W2").Value = 12
$ws.Range("X2").Value = 62175
$ws.Range("Y2").Value = 93
$ws.Range("Z2").Value = 865
$ws.Range("AA2").Value = 429
